$d = $word.ActiveDocument

$d.Content.Find.Execute("391×8=3128", $true, $false, $false, $false, $false, $true, 1, $false, "185×7=1295", 2) | Out-Null
$d.Content.Find.Execute("661×2=1322", $true, $false, $false, $false, $false, $true, 1, $false, "829×8=6632", 2) | Out-Null
$d.Content.Find.Execute("930×8=7440", $true, $false, $false, $false, $false, $true, 1, $false, "391×4=1564", 2) | Out-Null
$d.Content.Find.Execute("292×2=584", $true, $false, $false, $false, $false, $true, 1, $false, "973×7=6811", 2) | Out-Null
$d.Content.Find.Execute("568×8=4544", $true, $false, $false, $false, $false, $true, 1, $false, "832×9=7488", 2) | Out-Null
$d.Content.Find.Execute("688×3=2064", $true, $false, $false, $false, $false, $true, 1, $false, "177×7=1239", 2) | Out-Null
$d.Content.Find.Execute("862×4=3448", $true, $false, $false, $false, $false, $true, 1, $false, "903×9=8127", 2) | Out-Null
$d.Content.Find.Execute("568×4=2272", $true, $false, $false, $false, $false, $true, 1, $false, "408×2=816", 2) | Out-Null
$d.Content.Find.Execute("384×7=2688", $true, $false, $false, $false, $false, $true, 1, $false, "466×5=2330", 2) | Out-Null
$d.Content.Find.Execute("870×3=2610", $true, $false, $false, $false, $false, $true, 1, $false, "580×2=1160", 2) | Out-Null
$d.Content.Find.Execute("450×9=4050", $true, $false, $false, $false, $false, $true, 1, $false, "982×3=2946", 2) | Out-Null
$d.Content.Find.Execute("869×4=3476", $true, $false, $false, $false, $false, $true, 1, $false, "285×5=1425", 2) | Out-Null
$d.Content.Find.Execute("639×7=4473", $true, $false, $false, $false, $false, $true, 1, $false, "602×7=4214", 2) | Out-Null
$d.Content.Find.Execute("544×5=2720", $true, $false, $false, $false, $false, $true, 1, $false, "927×5=4635", 2) | Out-Null
$d.Content.Find.Execute("281×7=1967", $true, $false, $false, $false, $false, $true, 1, $false, "748×9=6732", 2) | Out-Null
$d.Content.Find.Execute("752×7=5264", $true, $false, $false, $false, $false, $true, 1, $false, "383×2=766", 2) | Out-Null
$d.Content.Find.Execute("330×2=660", $true, $false, $false, $false, $false, $true, 1, $false, "386×8=3088", 2) | Out-Null
$d.Content.Find.Execute("693×4=2772", $true, $false, $false, $false, $false, $true, 1, $false, "563×5=2815", 2) | Out-Null
$d.Content.Find.Execute("880×4=3520", $true, $false, $false, $false, $false, $true, 1, $false, "245×2=490", 2) | Out-Null
$d.Content.Find.Execute("948×8=7584", $true, $false, $false, $false, $false, $true, 1, $false, "792×4=3168", 2) | Out-Null
$d.Content.Find.Execute("741×5=3705", $true, $false, $false, $false, $false, $true, 1, $false, "675×7=4725", 2) | Out-Null
$d.Content.Find.Execute("287×4=1148", $true, $false, $false, $false, $false, $true, 1, $false, "211×7=1477", 2) | Out-Null
$d.Content.Find.Execute("123×2=246", $true, $false, $false, $false, $false, $true, 1, $false, "289×7=2023", 2) | Out-Null
$d.Content.Find.Execute("829×4=3316", $true, $false, $false, $false, $false, $true, 1, $false, "542×2=1084", 2) | Out-Null
$d.Content.Find.Execute("230×4=920", $true, $false, $false, $false, $false, $true, 1, $false, "374×4=1496", 2) | Out-Null
